$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be protected
# as Text so Excel does not coerce them to a numeric cell type -
# matches the source workbook, which stores these as literal text.
$textCells = @("D4", "D5", "D6", "D7", "D10", "D11", "D12", "D13", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.103.30"
$ws.Range("E2").Value = "  -5.43%  "
$ws.Range("D3").Value = "3.301.01"
$ws.Range("E3").Value = "  -7.89%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "181.94"
$ws.Range("E5").Value = "  -9.77%  "
$ws.Range("D6").Value = "518.59"
$ws.Range("E6").Value = "  -9.21%  "
$ws.Range("D7").Value = "0.589"
$ws.Range("E7").Value = "  -4.04%  "
$ws.Range("D8").Value = "3.299.90"
$ws.Range("E8").Value = "  -7.81%  "
$ws.Range("D10").Value = "0.615"
$ws.Range("E10").Value = "  -9.19%  "
$ws.Range("D11").Value = "58.22"
$ws.Range("E11").Value = "  -3.51%  "
$ws.Range("D12").Value = "0.131"
$ws.Range("E12").Value = "  -11.35%  "
$ws.Range("D13").Value = "0.0000251"
$ws.Range("E13").Value = "  -10.52%  "
$ws.Range("D14").Value = "9.02"
$ws.Range("E14").Value = "  -12.10%  "
$ws.Range("D15").Value = "3.818.82"
$ws.Range("E15").Value = "  -8.07%  "
$ws.Range("E16").Value = "  -4.47%  "
$ws.Range("D17").Value = "3.303.90"
$ws.Range("E17").Value = "  -7.68%  "
$ws.Range("D18").Value = "63.656.41"
$ws.Range("E18").Value = "  -5.74%  "
$ws.Range("D19").Value = "17.15"
$ws.Range("E19").Value = "  -10.16%  "
$ws.Range("D20").Value = "10.82"
$ws.Range("E20").Value = "  -11.40%  "
$ws.Range("D21").Value = "0.943"
$ws.Range("E21").Value = "  -11.07%  "
$ws.Range("D22").Value = "369.98"
$ws.Range("E22").Value = "  -8.23%  "
$ws.Range("D23").Value = "79.84"
$ws.Range("E23").Value = "  -5.64%  "
$ws.Range("D24").Value = "3.65"
$ws.Range("E24").Value = "  -13.16%  "
$ws.Range("D25").Value = "10.73"
$ws.Range("E25").Value = "  -15.96%  "
$ws.Range("E26").Value = "  -1.96%  "
$ws.Range("D27").Value = "3.75"
$ws.Range("E27").Value = "  -3.38%  "
$ws.Range("D28").Value = "2.63"
$ws.Range("E28").Value = "  -9.20%  "
$ws.Range("D29").Value = "11.12"
$ws.Range("E29").Value = "  -10.55%  "
$ws.Range("D30").Value = "8.27"
$ws.Range("E30").Value = "  -10.06%  "
$ws.Range("D31").Value = "645.37"
$ws.Range("E31").Value = "  -4.60%  "
$ws.Range("D32").Value = "28.39"
$ws.Range("E32").Value = "  -9.64%  "
$ws.Range("D33").Value = "6.66"
$ws.Range("E33").Value = "  -13.30%  "
$ws.Range("D34").Value = "11.04"
$ws.Range("E34").Value = "  -8.84%  "
$ws.Range("D35").Value = "59.21"
$ws.Range("E35").Value = "  -6.41%  "
$ws.Range("D36").Value = "0.103"
$ws.Range("E36").Value = "  -9.38%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "35.89"
$ws.Range("E38").Value = "  -13.56%  "
$ws.Range("D39").Value = "0.372"
$ws.Range("E39").Value = "  -9.22%  "
$ws.Range("D40").Value = "0.995"
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("E41").Value = "  -8.85%  "
$ws.Range("D42").Value = "2.785.81"
$ws.Range("E42").Value = "  -12.55%  "
$ws.Range("D43").Value = "2.66"
$ws.Range("E43").Value = "  -17.05%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0617"
$ws.Range("E44").Value = "  -18.89%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "2.58"
$ws.Range("E45").Value = "  -8.17%  "
$ws.Range("D46").Value = "0.0384"
$ws.Range("E46").Value = "  -6.48%  "
$ws.Range("D47").Value = "2.29"
$ws.Range("E47").Value = "  -15.25%  "
$ws.Range("D48").Value = "0.123"
$ws.Range("E48").Value = "  -5.66%  "
$ws.Range("D49").Value = "133.96"
$ws.Range("E49").Value = "  -3.49%  "
$ws.Range("D50").Value = "2.64"
$ws.Range("E50").Value = "  -2.59%  "
$ws.Range("D51").Value = "2.82"
$ws.Range("E51").Value = "  -8.15%  "

# Remove the temporary Text format so the cells end up back on the
# default (unstyled) cell format, same as the rest of the sheet.
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
